$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '44.315.25'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.91%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.372.29'
$ws.Range("D3").Style = "Normal"

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.694'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +6.29%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '243.96'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.53%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '76.43'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +6.30%  '

$ws.Range("E8").Value = '  +0.06%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.590'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +24.30%  '

$ws.Range("E10").Value = '  +4.68%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '58.07'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.17%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '31.90'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +16.41%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.44'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +17.40%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.108'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.32%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.728.07'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.02%  '

$ws.Range("E16").Value = '  +6.11%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.923'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +7.31%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.377.06'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.28%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '44.314.54'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.92%  '

$ws.Range("E20").Value = '  +2.69%  '

$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.70'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.22%  '

$ws.Range("B22").Value = 'Litecoin'
$ws.Range("C22").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '78.62'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.13%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '258.17'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.68%  '

$ws.Range("E24").Value = '  +0.09%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.58'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.75%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.72'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.60%  '

$ws.Range("E27").Value = '  +7.53%  '

$ws.Range("B28").Value = 'ImmutableX'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.73'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +13.17%  '

$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.31'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.79%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '23.14'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.71%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '175.42'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.33%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.130'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.30%  '

$ws.Range("E33").Value = '  +6.41%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.44'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +8.04%  '

$ws.Range("E35").Value = '  +8.93%  '

$ws.Range("E36").Value = '  +5.53%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.87'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.89%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.50'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.49%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.61'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.48%  '

$ws.Range("E40").Value = '  +7.77%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '19.17'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.64%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.12'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.05%  '

$ws.Range("E43").Value = '  +0.08%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.195'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +15.94%  '

$ws.Range("B45").Value = 'ARBITRUM'
$ws.Range("C45").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.23'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.70%  '

$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.101'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.37%  '

$ws.Range("E47").Value = '  +3.47%  '

$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.51'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +11.06%  '

$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '102.00'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.02%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.46'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.43%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.470.22'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.45%  '
